$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 112.84340945255747
$ws.Range("C4").Value = 52263.1431262652
$ws.Range("C5").Value = 49919.610559671804
$ws.Range("C6").Value = 2467.1431262652072
$ws.Range("C7").Value = 2343.7859699519468
$ws.Range("C12").Value = 0.5308083423660894
$ws.Range("C13").Value = 96.57603827332888
$ws.Range("C14").Value = 142.887868983031
$ws.Range("C15").Value = 95.44884708572306
$ws.Range("C16").Value = 22.017503337030206
$ws.Range("C17").Value = 100.0
$ws.Range("C18").Value = 38.39482262187991
$ws.Range("C21").Value = 0.5573494361384519
$ws.Range("C23").Value = 0.4023727896790188
$ws.Range("C24").Value = 16.40161372463008
$ws.Range("C25").Value = 18.64913214034498
$ws.Range("C27").Value = 4.476945800608419
$ws.Range("C28").Value = 13.595467385791132
$ws.Range("C32").Value = 0.3512109448371987
$ws.Range("C34").Value = 37.33126781772375
$ws.Range("C35").Value = 788.8369109468399
$ws.Range("C36").Value = 537.6088759926045
$ws.Range("C37").Value = 1.5680932681558488
$ws.Range("C38").Value = 270.00861556829545
$ws.Range("C39").Value = 564.3780993874409
$ws.Range("C40").Value = 21.583618213096315
$ws.Range("C41").Value = 122.50075931938022
$ws.Range("C42").Value = -0.03027056158975938
$ws.Range("C45").Value = 52226.065261806056
$ws.Range("C46").Value = 51437.228350859215
$ws.Range("C47").Value = 50899.6194748666
$ws.Range("C48").Value = 50898.051381598445
$ws.Range("C49").Value = 50628.04276603015
$ws.Range("C50").Value = 50063.664666642704
$ws.Range("C51").Value = 50042.08104842961
$ws.Range("C52").Value = 49919.58028911022
$ws.Range("C53").Value = 49919.61055967181
$ws.Range("C54").Value = 49919.610559671804
$ws.Range("C58").Value = 147.06581114238537
$ws.Range("C60").Value = 0.2223556780512614
$ws.Range("C62").Value = 2.126695251478622
$ws.Range("C64").Value = 0.09267083118915845
$ws.Range("C66").Value = 22.94891741218594
$ws.Range("C67").Value = 46611.9489095388
$ws.Range("C68").Value = 37403.91950088453
$ws.Range("C70").Value = 7022.615112008484
$ws.Range("C77").Value = 0.41383800477273924
$ws.Range("C78").Value = 1.0876918253545205
$ws.Range("C79").Value = 0.029563764164287156
$ws.Range("C80").Value = 0.07332278232273193
$ws.Range("C81").Value = 13.9981499809369
$ws.Range("C82").Value = 14.83429557496904
$ws.Range("C83").Value = 24073.62015582046
$ws.Range("C84").Value = 11059.118118010332
$ws.Range("C85").Value = 8221.467394931173
$ws.Range("C86").Value = 7807.31702248335
$ws.Range("C89").Value = 459.08300177402543
$ws.Range("C90").Value = 460.52731432613206
$ws.Range("C91").Value = 0.7788563201937818
$ws.Range("C92").Value = 0.7813066656763994
$ws.Range("C93").Value = 0.348875229484687
$ws.Range("C94").Value = 0.34125228922845885
$ws.Range("C95").Value = 0.025547831227245296
$ws.Range("C96").Value = 0.025369749949183916
$ws.Range("C97").Value = 13.655766956556047
$ws.Range("C98").Value = 13.45114910127193
$ws.Range("C99").Value = 8304.160256122626
$ws.Range("C100").Value = 8298.244601761078
$ws.Range("C101").Value = 8304.160256122626
$ws.Range("C102").Value = 8298.244601761078
$ws.Range("C109").Value = 0.5628463321985737
$ws.Range("C110").Value = 0.5625034934743327
$ws.Range("C111").Value = 0.03247844996867854
$ws.Range("C112").Value = 0.03246436440058862
$ws.Range("C113").Value = 17.329839716531104
$ws.Range("C114").Value = 17.326798286681804
$ws.Range("C115").Value = -311.0648126141039
$ws.Range("C116").Value = 809.9170717580963
$ws.Range("C117").Value = 7120.31593457726
$ws.Range("C118").Value = 7117.22002628589
$ws.Range("C125").Value = 0.40331486984889353
$ws.Range("C126").Value = 0.6098275704612877
$ws.Range("C127").Value = 0.029169799947729794
$ws.Range("C128").Value = 0.03814301431355167
$ws.Range("C129").Value = 13.826453063497352
$ws.Range("C130").Value = 15.987922859170169
$ws.Range("C133").Value = 7946.904670725951
$ws.Range("C134").Value = 6968.1702634312605
$ws.Range("C137").Value = 441.21068121021403
$ws.Range("C138").Value = 441.21068121021403
$ws.Range("C139").Value = 0.704285065228287
$ws.Range("C140").Value = 0.704285065228287
$ws.Range("C141").Value = 0.22132046496662092
$ws.Range("C142").Value = 0.217171901093037
$ws.Range("C143").Value = 0.022849750062349312
$ws.Range("C144").Value = 0.022772084157704518
$ws.Range("C145").Value = 9.685903100152585
$ws.Range("C146").Value = 9.536759990392046
$ws.Range("C147").Value = 11523.52312714743
$ws.Range("C148").Value = 11507.889093627313
$ws.Range("C149").Value = 11523.52312714743
$ws.Range("C150").Value = 11507.889093627313
$ws.Range("C157").Value = 0.5532821109211329
$ws.Range("C158").Value = 0.5525422808074715
$ws.Range("C159").Value = 0.03209357219044354
$ws.Range("C160").Value = 0.03206645570595816
$ws.Range("C161").Value = 17.23965495763301
$ws.Range("C162").Value = 17.231161618675724
$ws.Range("C163").Value = 1131.8962165037506
$ws.Range("C164").Value = 445.58197273788926
$ws.Range("C165").Value = 7033.948934994985
$ws.Range("C166").Value = 7027.268108060516
$ws.Range("C173").Value = 0.7188080623029574
$ws.Range("C174").Value = 0.717487619075011
$ws.Range("C175").Value = 0.0394401926933973
$ws.Range("C176").Value = 0.03937423095146737
$ws.Range("C177").Value = 18.22526750543218
$ws.Range("C178").Value = 18.22226369219466
$ws.Range("C179").Value = 6053.34895034462
$ws.Range("C180").Value = 6043.225028172676
$ws.Range("C181").Value = 6053.34895034462
$ws.Range("C182").Value = 6043.225028172676
$ws.Range("C189").Value = 0.5512790648823496
$ws.Range("C190").Value = 0.5512413716588563
$ws.Range("C191").Value = 0.03202015592117955
$ws.Range("C192").Value = 0.03201877437737485
$ws.Range("C193").Value = 17.21662649736534
$ws.Range("C194").Value = 17.216192136585192
$ws.Range("C195").Value = -58.59660992986095
$ws.Range("C196").Value = 13.594577309080737
$ws.Range("C197").Value = 7015.860993344899
$ws.Range("C198").Value = 7015.520615333634
$ws.Range("C201").Value = 117.3037566830413
$ws.Range("C203").Value = 0.17736618458277903
$ws.Range("C213").Value = 4558.735009648629
